$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.798.49'
$ws.Range('E2').Value = '  -0.96%  '
$ws.Range('D3').Value = '3.482.15'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '605.03'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.06'
$ws.Range('E6').Value = '  -2.14%  '
$ws.Range('D7').Value = '3.478.23'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.477'
$ws.Range('E9').Value = '  -1.62%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.141'
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.88'
$ws.Range('E11').Value = '  +4.83%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.417'
$ws.Range('E12').Value = '  -2.28%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '4.098.59'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('B14').Value = 'ShibaInu'
$ws.Range('C14').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000213'
$ws.Range('E14').Value = '  -0.90%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '31.14'
$ws.Range('E15').Value = '  -2.30%  '
$ws.Range('D16').Value = '3.488.82'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').Value = '66.901.10'
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.71'
$ws.Range('E19').Value = '  +8.30%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.27'
$ws.Range('E20').Value = '  -3.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.24'
$ws.Range('E21').Value = '  -1.29%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '431.32'
$ws.Range('E22').Value = '  -3.78%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.602'
$ws.Range('E23').Value = '  -3.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.70'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.623.21'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000117'
$ws.Range('E27').Value = '  -4.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.73'
$ws.Range('E28').Value = '  -1.33%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.12'
$ws.Range('E29').Value = '  -6.01%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.49'
$ws.Range('E30').Value = '  -0.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.56'
$ws.Range('E31').Value = '  -5.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.26%  '
$ws.Range('E33').Value = '  -2.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.30'
$ws.Range('E34').Value = '  -0.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.78'
$ws.Range('E35').Value = '  -2.72%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.76'
$ws.Range('E36').Value = '  -6.71%  '
$ws.Range('B37').Value = 'USDe'
$ws.Range('C37').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').Value = '  +0.00%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.93'
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('E39').Value = '  +0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '175.33'
$ws.Range('E40').Value = '  -1.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0890'
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.34'
$ws.Range('E42').Value = '  -0.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.01'
$ws.Range('E43').Value = '  -11.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.891'
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '46.35'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '27.61'
$ws.Range('E46').Value = '  -10.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.21'
$ws.Range('E47').Value = '  -6.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.34'
$ws.Range('E48').Value = '  -3.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.40'
$ws.Range('E49').Value = '  -3.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.977'
$ws.Range('E50').Value = '  -1.52%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.244'
$ws.Range('E51').Value = '  -1.97%  '
